$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.455.94"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "'1.565.73"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'289.26"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.3688"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("D8").Value = "'50.20"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "'0.3379"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "'1.143"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'21.12"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'6.009"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "'6.974"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "'1.571.91"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'0.00001115"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "'90.17"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "'0.06752"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "'6.374"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").Value = "'16.30"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "'12.08"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "'22.447.73"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'2.398"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").Value = "'2.635"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").Value = "'19.89"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'149.21"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").Value = "'5.060"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "'124.61"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'1.744.70"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'1.054"
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("D33").Value = "'6.182"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").Value = "'9.718"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "'0.08333"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'0.02471"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "'0.2284"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("D39").Value = "'1.338"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("D40").Value = "'0.06474"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").Value = "'5.392"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "'11.21"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "'0.6187"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'13.89"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "'3.769"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'0.5814"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "'125.32"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "'1.227"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'0.07340"
$ws.Range("E51").Value = "  +0.61%  "
